# Add a new row (14) to Sheet1 with a new query in column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New query text added as cell C14 (row 14), matching the style (vertical
# center) used by the rest of column C.
$ws.Range("C14").Value = "What lyrics come after this line ""Loving him was blue like I'd never known""?"
$ws.Range("C14").VerticalAlignment = -4108

# Mirror the author's new selection state: the active cell in the
# (frozen) bottom-right pane becomes C13.
$ws.Range("C13").Select()
